# fix errors in data, visualize clean data
#
# Rows 57-92 on the "dipole" sheet had columns B (param 1) and C (param 2)
# swapped relative to the correct pattern used everywhere else in the table
# (B cycles 2,5,10,50,150 while C holds the group's constant). Swap B<->C
# back for every one of these rows to fix the bad data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 57; $r -le 92; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $bVal = $bCell.Value2
    $cVal = $cCell.Value2
    $bCell.Value = $cVal
    $cCell.Value = $bVal
}

# Scroll the view down and select F88 so the newly-cleaned data is visible.
$win = $excel.ActiveWindow
$win.ScrollRow = 70
$win.ScrollColumn = 1
$ws.Range("F88").Select() | Out-Null
